# Update BOM sheet: designator range R29-R32 -> R29-R33 (n_bits resistor
# no longer generated, so the group now spans one more designator), and
# swap the stocked/alternate-part quantities for the two TAR5S1xUTE85LF
# regulator rows (U6 now sources the -18 part instead of the -16 part).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Row 6/7: swap QTY between the primary and alternate part for U6
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 1

# Row 27: designator text + quantity for the R29-R33 resistor group
$ws.Range("A27").Value = "R29-R33"
$ws.Range("B27").Value = 5

# Restore the last-used selection seen in the saved workbook
$ws.Range("D12").Select()
